# Apply the PrivateCar_example.xlsx edits:
#  - add two new columns (U: P_discharge, V: P_discharge_variance) with header + data
#  - update several existing numeric columns (A,B,C,D,E,N,O,Q,R) and the T text column
#    for every data row (rows 2-11)
#  - move the active selection to V19 to mirror the saved view state

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column headers (U1, V1) ---
$ws.Range("U1").Value = "P_discharge"
$ws.Range("V1").Value = "P_discharge_variance"

# Text used for column T (shared strings already present in the workbook)
$T18 = "17,18,19,20,21"
$T19 = "25,26,27,28,29"

# --- Per-row data (row number -> hashtable of column letter/value) ---
$rows = @{
    2  = @{ A = 1;  B = 60;   C = 1; D = 7.2; E = 0.1; N = 60;   O = 1; Q = 50; R = 3; T = $T18; U = 60;   V = 1 }
    3  = @{ A = 2;  B = 61.9; C = 1; D = 9;   E = 0.1; N = 51.6; O = 1; Q = 50; R = 3; T = $T19; U = 51.6; V = 1 }
    4  = @{ A = 3;  B = 76.8; C = 1; D = 7.8; E = 0.1; N = 76.8; O = 1; Q = 50; R = 3; T = $T18; U = 76.8; V = 1 }
    5  = @{ A = 4;  B = 60;   C = 1; D = 7.2; E = 0.1; N = 60;   O = 1; Q = 50; R = 3; T = $T18; U = 60;   V = 1 }
    6  = @{ A = 5;  B = 47.5; C = 1; D = 7.2; E = 0.1; N = 95;   O = 1; Q = 50; R = 3; T = $T18; U = 95;   V = 1 }
    7  = @{ A = 6;  B = 60;   C = 1; D = 7.2; E = 0.1; N = 60;   O = 1; Q = 50; R = 3; T = $T18; U = 60;   V = 1 }
    8  = @{ A = 7;  B = 61.9; C = 1; D = 9;   E = 0.1; N = 51.6; O = 1; Q = 50; R = 3; T = $T19; U = 51.6; V = 1 }
    9  = @{ A = 8;  B = 76.8; C = 1; D = 7.8; E = 0.1; N = 76.8; O = 1; Q = 50; R = 3; T = $T18; U = 76.8; V = 1 }
    10 = @{ A = 9;  B = 61.9; C = 1; D = 9;   E = 0.1; N = 51.6; O = 1; Q = 50; R = 3; T = $T19; U = 51.6; V = 1 }
    11 = @{ A = 10; B = 47.5; C = 1; D = 7.2; E = 0.1; N = 95;   O = 1; Q = 50; R = 3; T = $T18; U = 95;   V = 1 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    foreach ($col in $vals.Keys) {
        $ws.Range("$col$r").Value = $vals[$col]
    }
}

# --- Column width for the new column V (matches col min=22 max=22 width=12) ---
# (ColumnWidth uses character units that get re-quantized on save; 79/7 round-trips to width=12)
$ws.Columns.Item(22).ColumnWidth = 79/7

# --- Restore the saved selection state (activeCell V19) ---
$ws.Range("V19").Select() | Out-Null
